# Update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.016.76"
Set-TextValue "D3" "1.644.67"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.18%  "
Set-TextValue "D5" "215.02"
$ws.Range("E5").Value = "  +2.08%  "
Set-TextValue "D6" "0.5220"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.19%  "
Set-TextValue "D8" "0.2608"
$ws.Range("E8").Value = "  -0.24%  "
Set-TextValue "D9" "0.06353"
$ws.Range("E9").Value = "  +0.34%  "
Set-TextValue "D10" "20.79"
$ws.Range("E10").Value = "  -1.71%  "
Set-TextValue "D11" "0.07660"
$ws.Range("E11").Value = "  +1.34%  "
Set-TextValue "D12" "1.645.25"
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("E13").Value = "  -0.07%  "
Set-TextValue "D14" "1.867.51"
Set-TextValue "D15" "0.5533"
$ws.Range("E15").Value = "  +1.63%  "
Set-TextValue "D16" "0.0₅8303"
$ws.Range("E16").Value = "  +3.29%  "
Set-TextValue "D17" "64.78"
$ws.Range("E17").Value = "  -2.56%  "
Set-TextValue "D18" "26.033.07"
$ws.Range("E18").Value = "  -0.49%  "
Set-TextValue "D20" "4.713"
$ws.Range("E20").Value = "  -0.67%  "
Set-TextValue "D21" "188.16"
$ws.Range("E21").Value = "  +0.30%  "
Set-TextValue "D22" "10.17"
$ws.Range("E22").Value = "  -1.02%  "
Set-TextValue "D23" "6.246"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -0.19%  "
Set-TextValue "D25" "145.18"
$ws.Range("E25").Value = "  -2.87%  "
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("E29").Value = "  +1.28%  "
Set-TextValue "D30" "0.05955"
Set-TextValue "D31" "1.265"
$ws.Range("E31").Value = "  -1.36%  "
Set-TextValue "B32" "Filecoin"
Set-TextValue "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "3.395"
$ws.Range("E32").Value = "  -0.87%  "
Set-TextValue "B33" "InternetComputer(DFINITY)"
Set-TextValue "C33" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "3.403"
$ws.Range("E33").Value = "  -2.97%  "
Set-TextValue "D34" "1.650"
$ws.Range("E34").Value = "  +0.07%  "
Set-TextValue "D35" "0.9936"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("E36").Value = "  -0.21%  "
Set-TextValue "D37" "2.752"
$ws.Range("E37").Value = "  -0.38%  "
Set-TextValue "D38" "0.5621"
$ws.Range("E38").Value = "  -6.55%  "
$ws.Range("E39").Value = "  -0.36%  "
Set-TextValue "D40" "5.849"
$ws.Range("E40").Value = "  -3.52%  "
Set-TextValue "D41" "0.8530"
$ws.Range("E41").Value = "  -1.27%  "
Set-TextValue "D43" "1.025.17"
$ws.Range("E43").Value = "  -8.15%  "
Set-TextValue "D44" "98.49"
$ws.Range("E44").Value = "  -2.12%  "
Set-TextValue "D45" "1.794.61"
$ws.Range("E45").Value = "  -1.50%  "
Set-TextValue "D46" "0.0₈110"
$ws.Range("E46").Value = "  +0.84%  "
Set-TextValue "D47" "55.67"
$ws.Range("E47").Value = "  +0.23%  "
Set-TextValue "D49" "8.090"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("E50").Value = "  -2.03%  "
Set-TextValue "D51" "0.4217"
$ws.Range("E51").Value = "  -0.55%  "
